$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# optimization_parameters sheet (7th tab) gets restructured:
#  - header row loses its duplicated "value" cells in C1:F1
#  - "Model"/"Sigmoid" row is renamed to "production_function"/"Sigmoid"
#  - a new "L_curve" / 1 row is inserted right below it
#  - the old "Deletion" / 0 / 3 row near the bottom is removed
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("optimization_parameters")

# Drop the extra duplicated "value" cells from the header row.
$ws.Range("C1:F1").ClearContents()

# Remove the obsolete "Deletion" row entirely (row 16: Deletion / 0 / 3).
$ws.Rows("16:16").Delete()

# Rename the "Model" row to "production_function" (value text stays "Sigmoid").
$ws.Range("A8").Value = "production_function"

# Insert the new "L_curve" row right after it, with value 1.
$ws.Rows("9:9").Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 1
$ws.Range("B9").NumberFormat = "0.00E+00"

# This sheet becomes the active tab/selected cell.
$ws.Activate()
$ws.Range("B10").Select()
